$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row: "..._old" -> "..._FV2304", "..._new" -> "..._FV2310"
#    (column K / "diff" is left untouched)
# ---------------------------------------------------------------------------
$headerMap = @{
    "A1" = "Segmentname_FV2304"
    "B1" = "Segmentgruppe_FV2304"
    "C1" = "Segment_FV2304"
    "D1" = "Datenelement_FV2304"
    "E1" = "Segment ID_FV2304"
    "F1" = "Code_FV2304"
    "G1" = "Qualifier_FV2304"
    "H1" = "Beschreibung_FV2304"
    "I1" = "Bedingungsausdruck_FV2304"
    "J1" = "Bedingung_FV2304"
    "L1" = "Segmentname_FV2310"
    "M1" = "Segmentgruppe_FV2310"
    "N1" = "Segment_FV2310"
    "O1" = "Datenelement_FV2310"
    "P1" = "Segment ID_FV2310"
    "Q1" = "Code_FV2310"
    "R1" = "Qualifier_FV2310"
    "S1" = "Beschreibung_FV2310"
    "T1" = "Bedingungsausdruck_FV2310"
    "U1" = "Bedingung_FV2310"
}
foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# ---------------------------------------------------------------------------
# 2) Freeze the header row (split below row 1, top-left of the scrollable
#    area is A2) and keep the active pane selection on A1.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

# ---------------------------------------------------------------------------
# 3) Turn the used range A1:U79 into a native Excel Table ("Table1") so the
#    data can be filtered/sorted, matching the regenerated AHB workbook.
# ---------------------------------------------------------------------------
$dataRange = $ws.Range("A1:U79")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
